# Add drag polars for takeoff and landing (flap-deflection-angle column, A)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- "plain" flap block (rows 2-6): unstyled numeric flap angles ---
$ws.Range("A2").Value = 30
$ws.Range("A3").Value = 35
$ws.Range("A4").Value = 40
$ws.Range("A5").Value = 45
$ws.Range("A6").Value = 50

# --- "split" flap block: the "split" label moves up to its own row (9),
#     numeric flap angles (black, non-theme font) take over rows 10-14 ---
$ws.Range("A9").Value = "split"

$ws.Range("A10").Value = 30
$ws.Range("A11").Value = 35
$ws.Range("A12").Value = 40
$ws.Range("A13").Value = 45
$ws.Range("A14").Value = 50

# --- "slotted" flap block (rows 18-22): black, non-theme font flap angles ---
$ws.Range("A18").Value = 30
$ws.Range("A19").Value = 35
$ws.Range("A20").Value = 40
$ws.Range("A21").Value = 45
$ws.Range("A22").Value = 50

# Apply the explicit black font colour to the styled flap-angle cells
# (matches the extra font + cellXfs entry introduced in styles.xml).
# Done as two single-area calls since a multi-area union range only
# formats its first area.
$ws.Range("A10:A14").Font.Color = 0
$ws.Range("A18:A22").Font.Color = 0

# --- sheet view tweaks ---
$ws.Range("E1").Select()
$excel.ActiveWindow.Zoom = 150
